# Update the "想去人数" (want-to-go count) figures in column F for both the
# "展览" sheet and the "全部类型" sheet (which mirrors the same data).
# Mapping of row -> new value:
#   F2:  8320 -> 8331
#   F3:  7750 -> 7762
#   F4:  120  -> 122
#   F10: 161  -> 162
#   F13: 125  -> 127
#   F14: 1311 -> 1317
#   F16: 50   -> 51

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 8331
    3  = 7762
    4  = 122
    10 = 162
    13 = 127
    14 = 1317
    16 = 51
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
